$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 94: date 2025-10-17 (serial 45947), station "四方坪站充电量(kw)"
$ws.Range("A94").Value = 45947
$ws.Range("B94").Value = "四方坪站充电量(kw)"
$row94 = @(650.25500000000011, 930.06399999999974, 397.02599999999995, 435.1, 297.50799999999998, 679.91599999999983, 445.98900000000009, 260.89699999999999, 246.73599999999999, 155.21100000000001, 132.928, 58.849999999999994, 899.13500000000022, 1084.961, 425.56, 308.702, 195.36099999999999, 267.88300000000004, 118.155, 91.419999999999987, 69.08, 30.84, 102.51600000000001, 77.02)
$cols = @("C","D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z")
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "94").Value = $row94[$i]
}

# New row 95: date 2025-10-17 (serial 45947), station "高岭站充电量(kw)"
$ws.Range("A95").Value = 45947
$ws.Range("B95").Value = "高岭站充电量(kw)"
$row95 = @(283.77299999999997, 178.89300000000003, 86.751000000000005, 151.09300000000002, 105.629, 217.26299999999998, 283.42399999999998, 125.997, 277.68499999999995, 96.781000000000006, 108.54900000000001, 180.06399999999999, 473.57900000000001, 418.93400000000003, 127.012, 84.215999999999994, 173.839, 36.914000000000001, 33.704000000000001, 0, 11.127000000000001, 15.148999999999999, 0, 0)
for ($i = 0; $i -lt $cols.Length; $i++) {
    $ws.Range($cols[$i] + "95").Value = $row95[$i]
}

$ws.Range("I98").Select() | Out-Null
